# Update cryptocurrency price/volume data as scraped on Wed Mar 22 19:53:41 UTC 2023
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'27.301.08"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -3.31%  "

$ws.Range("D3").Value = "'1.748.58"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -2.80%  "

$ws.Range("D4").Value = "'1.019"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +1.56%  "

$ws.Range("D5").Value = "'323.77"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -4.32%  "

$ws.Range("D6").Value = "'1.013"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.45%  "

$ws.Range("D7").Value = "'0.4193"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -12.23%  "

$ws.Range("D8").Value = "'0.3579"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.71%  "

$ws.Range("D9").Value = "'44.81"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.45%  "

$ws.Range("D10").Value = "'1.109"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -2.98%  "

$ws.Range("D11").Value = "'0.07315"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -4.55%  "

$ws.Range("D12").Value = "'1.015"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.45%  "

$ws.Range("D13").Value = "'21.44"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -4.89%  "

$ws.Range("D14").Value = "'6.044"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -3.90%  "

$ws.Range("D15").Value = "'7.184"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.47%  "

$ws.Range("D16").Value = "'1.757.44"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -2.23%  "

$ws.Range("D17").Value = "'0.00001049"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -3.98%  "

$ws.Range("D18").Value = "'82.44"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.93%  "

$ws.Range("B19").Value = "Dai"
$ws.Range("C19").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D19").Value = "'1.012"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.31%  "

$ws.Range("B20").Value = "TRON"
$ws.Range("C20").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D20").Value = "'0.05948"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -11.40%  "

$ws.Range("D21").Value = "'16.67"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -3.60%  "

$ws.Range("D22").Value = "'6.045"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -5.50%  "

$ws.Range("D23").Value = "'27.388.45"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -3.03%  "

$ws.Range("D24").Value = "'11.14"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -7.13%  "

$ws.Range("D25").Value = "'2.404"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.27%  "

$ws.Range("B26").Value = "EthereumClassic"
$ws.Range("C26").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D26").Value = "'19.81"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -3.63%  "

$ws.Range("B27").Value = "Monero"
$ws.Range("C27").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D27").Value = "'149.84"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.79%  "

$ws.Range("B28").Value = "LidoDAOToken"
$ws.Range("C28").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D28").Value = "'2.326"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -3.32%  "

$ws.Range("D29").Value = "'1.959.08"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -2.28%  "

$ws.Range("D30").Value = "'1.247"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.99%  "

$ws.Range("D31").Value = "'125.71"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -5.97%  "

$ws.Range("D32").Value = "'3.700"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -8.43%  "

$ws.Range("D33").Value = "'0.08968"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -6.67%  "

$ws.Range("D34").Value = "'5.478"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -7.42%  "

$ws.Range("D35").Value = "'12.30"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +1.27%  "

$ws.Range("D36").Value = "'0.2138"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.46%  "

$ws.Range("D37").Value = "'0.06075"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -3.46%  "

$ws.Range("B38").Value = "VeChain"
$ws.Range("C38").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D38").Value = "'0.02244"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -5.95%  "

$ws.Range("B39").Value = "TheSandbox"
$ws.Range("C39").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D39").Value = "'0.6382"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -4.19%  "

$ws.Range("D40").Value = "'4.959"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -4.69%  "

$ws.Range("D41").Value = "'1.427"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -3.64%  "

$ws.Range("D42").Value = "'1.010"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.16%  "

$ws.Range("D43").Value = "'1.166"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -4.00%  "

$ws.Range("D44").Value = "'7.921"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.88%  "

$ws.Range("D45").Value = "'13.54"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -3.91%  "

$ws.Range("D46").Value = "'3.769"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -2.51%  "

$ws.Range("D47").Value = "'0.5827"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -4.84%  "

$ws.Range("D48").Value = "'123.05"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -4.06%  "

$ws.Range("D49").Value = "'1.924"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -5.53%  "

$ws.Range("D50").Value = "'0.06852"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -3.52%  "

$ws.Range("D51").Value = "'1.092"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -6.80%  "
